$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Reorder GUILHERME (row 10) / GILSON (row 11) and update
#    GUILHERME's balance from 25376.92 to 22000.
#    Result: GILSON (23860.68) now comes first (row 10),
#            GUILHERME (22000) comes second (row 11) - keeps the
#            sheet sorted descending by balance.
# -----------------------------------------------------------------
$ws.Rows.Item(11).EntireRow.Delete()
$ws.Rows.Item(10).EntireRow.Insert()
$ws.Cells.Item(10, 1).NumberFormat = "@"
$ws.Cells.Item(10, 1).Value = "004474776"
$ws.Cells.Item(10, 2).Value = "GILSON"
$ws.Cells.Item(10, 3).Value = 23860.68

$ws.Cells.Item(11, 3).Value = 22000

# -----------------------------------------------------------------
# 2) Move DANIELE's row (originally row 293, balance 38.82) up so it
#    sits just above GABRIEL (row 33), and update her balance to
#    1538.82 - keeps the sheet sorted descending by balance.
# -----------------------------------------------------------------
$ws.Rows.Item(293).EntireRow.Delete()
$ws.Rows.Item(33).EntireRow.Insert()
$ws.Cells.Item(33, 1).NumberFormat = "@"
$ws.Cells.Item(33, 1).Value = "004398174"
$ws.Cells.Item(33, 2).Value = "DANIELE"
$ws.Cells.Item(33, 3).Value = 1538.82

# -----------------------------------------------------------------
# 3) Remove RICARDO (row 356) and MARCEL (row 357) entirely - their
#    row numbers are unchanged because steps above net to zero row
#    shift below row 33.
# -----------------------------------------------------------------
$ws.Rows.Item(357).EntireRow.Delete()
$ws.Rows.Item(356).EntireRow.Delete()
